$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "dox" in D1, copying the header formatting (style) used by C1
$ws.Range("D1").Value = "dox"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "dox"

# Row 2: split "dusp11 -dox foldchange" -> B2 "dusp11", D2 "-dox"
$ws.Range("B2").Value = "dusp11"
$ws.Range("D2").Value = "-dox"

# Row 3: split "ifnb -dox foldchange" -> B3 "ifnb", D3 "-dox"
$ws.Range("B3").Value = "ifnb"
$ws.Range("D3").Value = "-dox"

# Row 4: split "mx1 -dox foldchange" -> B4 "mx1", D4 "-dox"
$ws.Range("B4").Value = "mx1"
$ws.Range("D4").Value = "-dox"

# Row 5: split "dusp11 +dox foldchange" -> B5 "dusp11", D5 "+dox"
$ws.Range("B5").Value = "dusp11"
$ws.Range("D5").Value = "+dox"

# Row 6: split "ifnb +dox foldchange" -> B6 "ifnb", D6 "+dox"
$ws.Range("B6").Value = "ifnb"
$ws.Range("D6").Value = "+dox"

# Row 7: split "mx1 +dox foldchange" -> B7 "mx1", D7 "+dox"
$ws.Range("B7").Value = "mx1"
$ws.Range("D7").Value = "+dox"
